$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update Riders (column C) and Average (column D) values per the ridership run on 20161026.
$ws.Range("C2").Value = 172
$ws.Range("D2").Value = 101.75

$ws.Range("C3").Value = 258
$ws.Range("D3").Value = 106.89

$ws.Range("C4").Value = 240
$ws.Range("D4").Value = 112.06

$ws.Range("C5").Value = 220
$ws.Range("D5").Value = 110.3

$ws.Range("C6").Value = 212
$ws.Range("D6").Value = 107

$ws.Range("D7").Value = 49.79

$ws.Range("C8").Value = 74
$ws.Range("D8").Value = 38.81

$wb.Save()
